# V 2.0.2 - se arreglo la fecha y hora de reimpresion
# Update patient record fields on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient identification ---
$ws.Range("A6").Value = "ALVARADO "
$ws.Range("C6").Value = "FLORES"
$ws.Range("E6").Value = "JUAN"
$ws.Range("G6").Value = "CARLOS"
$ws.Range("I6").Value = "/201762647"

# --- Direccion actual ---
$ws.Range("A8").Value = "MZ N"
$ws.Range("D8").Value = "EL MILAGRO ZONA 10 "
$ws.Range("F8").Value = "MIXCO"
$ws.Range("H8").Value = ""

# --- Fecha / edad / lugar de nacimiento ---
$ws.Range("A12").Value = "1970-03-06"
$ws.Range("F12").Value = "47"
$ws.Range("H12").Value = "ESTA CAPITAL"

# --- Estado civil / ocupacion / cedula ---
$ws.Range("A14").Value = "Casado"
$ws.Range("D14").Value = "ALBAÑIL"
$ws.Range("H14").Value = "1969701480101"

# --- Nombre del padre / madre ---
$ws.Range("A18").Value = "ALVARADO"
$ws.Range("F18").Value = "TEODORA FLORES"

# --- Contacto de emergencia ---
$ws.Range("A20").Value = "TEODORA FLORES"
$ws.Range("J20").Value = "54138784"

# --- Fecha / hora de ingreso y servicio ---
$ws.Range("A24").Value = "24/10/2017"
$ws.Range("C24").Value = "15:36:23"
$ws.Range("D24").Value = ""
